{"js": "// Remove the \" \u2022 {{ job.location }}\" run sequence that follows\n// \"{{ job.company }}\" in the earlier-experience job line, leaving:\n// \"{{ job.company }} \u2022 {{ job.title }} \u2022 {{ job.dates }}\"\nconst body = context.document.body;\nconst results = body.search(\" \u2022 {{ job.location }}\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (const r of results.items) {\n  r.delete();\n}\nawait context.sync();\n", "ps1": "# Remove the \" \u2022 {{ job.location }}\" text (a leading space, the bullet\n# character, and the location placeholder) that follows \"{{ job.company }}\"\n# in the earlier-experience job line, leaving:\n# \"{{ job.company }} \u2022 {{ job.title }} \u2022 {{ job.dates }}\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \" \u2022 {{ job.location }}\"\n$find.Replacement.Text = \"\"\n\n$result = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
